# faktura_mal.xlsx - "div changes and add ons"
#
# Sheet "4" (the invoice line-item sheet) gets:
#  - a new "Totalt_antall" (count total) column in M, mirroring the
#    existing "Totalt" (money total) column in L
#  - B2 turned into a lookup-style formula instead of a hard-coded name
#  - H2/H3 given their Momspris (VAT amount) formula F*G
#  - L2's total formula switched from summing I2:I3 to summing F2:F3
#  - I3's running-total formula switched from SUM(I2) to F3+H3
#  - J3's date corrected back to the same date as J2
#  - two more rows (4 and 5) filled in with elevnummer/name so the
#    third student entry has its own row
#  - the L2 currency style gets a plain "kr" #,##0.00 number format
#  - the active selection left on B4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4")

# --- header row: new "Totalt_antall" column -------------------------------
$ws.Range("M1").Value = "Totalt_antall"

# --- row 2 ------------------------------------------------------------
$ws.Range("B2").Formula = '=IF(A2=1, "Anette Evensen", "N/A")'
$ws.Range("H2").Formula = "=F2*G2"
$ws.Range("L2").Formula = "=SUM(F2:F3)"
$ws.Range("M2").Formula = "=SUM(C2:C3)"

# --- row 3 ------------------------------------------------------------
$ws.Range("H3").Formula = "=F3*G3"
$ws.Range("I3").Formula = "=F3+H3"
$ws.Range("J3").Value = 43937

# --- row 4: new entry (same student as rows 2/3) -----------------------
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Anette Evensen"

# --- row 5: elevnummer filled in ---------------------------------------
$ws.Range("A5").Value = 1

# --- number format for the "Totalt" cell --------------------------------
$ws.Range("L2").NumberFormat = """kr"" #,##0.00"

# --- leave the selection on B4, matching the saved file -----------------
$ws.Range("B4").Select() | Out-Null
